$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 for FY_4.png, shifting existing rows 36-44 down to 37-45
$ws.Rows.Item(36).Insert()

$data = @(
  @('AK_1.png', 0.442, 0.001, 0.795, 0.533, 'Akhlak Kamiswara', 'Benar'),
  @('AK_2.png', 0.859, 0.003, 0.853, 0.867, 'Akhlak Kamiswara', 'Benar'),
  @('AK_3.png', 0.835, 0.003, 0.856, 0.933, 'Akhlak Kamiswara', 'Benar'),
  @('AK_4.png', 0.438, 0.001, 0.833, 0.533, 'Akhlak Kamiswara', 'Benar'),
  @('AK_5.png', 0.453, 0.001, 0.854, 0.6, 'Akhlak Kamiswara', 'Benar'),
  @('MIB_1.png', 1.669, 0.005, 0.658, 0.867, 'Muhammad Iqbal Baqi', 'Benar'),
  @('MIB_2.png', 1.191, 0.004, 0.742, 0.867, 'Muhammad Iqbal Baqi', 'Benar'),
  @('MIB_3.png', 1.754, 0.006, 0.95, 1, 'Muhammad Iqbal Baqi', 'Benar'),
  @('MIB_4.png', 1.157, 0.004, 0.746, 0.867, 'Muhammad Iqbal Baqi', 'Benar'),
  @('MIB_5.png', 1.626, 0.005, 0.718, 0.867, 'Muhammad Iqbal Baqi', 'Benar'),
  @('AAH_1.png', 0.957, 0.003, 0.773, 0.933, 'Andrea Ayunove Hutami', 'Benar'),
  @('AAH_2.png', 1.098, 0.004, 0.932, 1, 'Andrea Ayunove Hutami', 'Benar'),
  @('AAH_3.png', 0.886, 0.003, 0.787, 1, 'Andrea Ayunove Hutami', 'Benar'),
  @('TI_1.png', 0.923, 0.003, 0.792, 0.6, 'Toni Ismail', 'Benar'),
  @('TI_2.png', 0.947, 0.003, 0.831, 0.733, 'Toni Ismail', 'Benar'),
  @('TI_3.png', 0.708, 0.002, 0.934, 0.6, 'Toni Ismail', 'Benar'),
  @('TI_4.png', 0.692, 0.002, 0.822, 0.8, 'Toni Ismail', 'Benar'),
  @('TI_5.png', 1.032, 0.003, 0.87, 1, 'Toni Ismail', 'Benar'),
  @('RAS_1.png', 0.646, 0.002, 0.844, 0.4, 'Tidak Diketahui', 'Salah'),
  @('RAS_2.png', 1.047, 0.004, 0.74, 0.933, 'Ridha Ayu Salsabila', 'Benar'),
  @('RAS_3.png', 0.518, 0.002, 0.781, 0.467, 'Tidak Diketahui', 'Salah'),
  @('RAS_4.png', 1.287, 0.004, 0.693, 0.4, 'Tidak Diketahui', 'Salah'),
  @('RAS_5.png', 1.168, 0.004, 0.78, 1, 'Ridha Ayu Salsabila', 'Benar'),
  @('RR_1.png', 1.509, 0.005, 0.864, 0.867, 'Rafiqo Rapitasari', 'Benar'),
  @('RR_2.png', 1.348, 0.004, 0.888, 0.867, 'Rafiqo Rapitasari', 'Benar'),
  @('RR_3.png', 1.106, 0.004, 0.517, 0.733, 'Rafiqo Rapitasari', 'Benar'),
  @('RR_4.png', 1.391, 0.005, 0.875, 0.933, 'Rafiqo Rapitasari', 'Benar'),
  @('RR_5.png', 1.449, 0.005, 0.875, 0.933, 'Rafiqo Rapitasari', 'Benar'),
  @('AR_1.png', 0.858, 0.003, 0.799, 1, 'Arizli Romadhon', 'Benar'),
  @('GA_1.png', 1.655, 0.005, 0.881, 1, 'Gege Ardiyansyah', 'Benar'),
  @('GA_2.png', 0.691, 0.002, 0.783, 0.333, 'Tidak Diketahui', 'Salah'),
  @('GA_3.png', 0.73, 0.002, 0.761, 0.267, 'Tidak Diketahui', 'Salah'),
  @('FY_1.png', 1.127, 0.004, 0.834, 0.333, 'Tidak Diketahui', 'Salah'),
  @('FY_2.png', 1.681, 0.006, 0.857, 0.467, 'Tidak Diketahui', 'Salah'),
  @('FY_3.png', 1.499, 0.005, 0.873, 0.533, 'Fanny Yusuf', 'Benar'),
  @('FY_4.png', 1.306, 0.004, 0.802, 0.467, 'Tidak Diketahui', 'Salah'),
  @('TO_1.png', 0.779, 0.003, 0.792, 0.8, 'Tiara Oktavian', 'Benar'),
  @('TO_2.png', 1.048, 0.003, 0.817, 1, 'Tiara Oktavian', 'Benar'),
  @('TO_3.png', 0.847, 0.003, 0.815, 1, 'Tiara Oktavian', 'Benar'),
  @('TO_4.png', 7.245, 0.024, 0.484, 0.4, 'Tidak Diketahui', 'Salah'),
  @('TO_5.png', 7.107, 0.023, 0.493, 0.333, 'Tidak Diketahui', 'Salah'),
  @('TD_1.png', 3.783, 0.013, 0.451, 0.333, 'Tidak Diketahui', 'Benar'),
  @('TD_2.png', 5.006, 0.017, 0.479, 0.333, 'Tidak Diketahui', 'Benar'),
  @('TD_3.png', 1.088, 0.004, 0.716, 0.267, 'Tidak Diketahui', 'Benar'),
  @('TD_4.png', 0.992, 0.003, 0.7, 0.267, 'Tidak Diketahui', 'Benar'),
)
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 1
  $values = $data[$i]
  for ($j = 0; $j -lt $values.Count; $j++) {
    $ws.Cells.Item($row, $j + 1).Value = $values[$j]
  }
}
